$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

# New weekly data rows to append (Player_1, Points_1, Player_2, Points_2)
$data = @(
    @(3, 6, 5, 14),
    @(5, 4, 4, 16),
    @(4, 8, 3, 12),
    @(2, 16, 3, 4),
    @(4, 7, 5, 13),
    @(7, 14, 6, 6),
    @(5, 19, 4, 1),
    @(4, 7, 5, 13),
    @(2, 6, 3, 14),
    @(1, 14, 4, 6),
    @(7, 5, 5, 15),
    @(3, 18, 6, 2),
    @(4, 5, 3, 15),
    @(6, 4, 4, 16)
)

$startRow = 1005
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$lastRow = $startRow + $data.Count - 1

# Scroll/selection state matching the saved view
$ws.Range("G1014").Select()
$excel.ActiveWindow.ScrollRow = 1004
